$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.00914180373853
$ws.Range("D2").Value = 10.36821571452253
$ws.Range("E2").Value = 17.5518400219533
$ws.Range("F2").Value = 32.94890039159539
$ws.Range("G2").Value = 3.65284475495932
$ws.Range("J2").Value = 12.35659796914091
$ws.Range("L2").Value = 11.27771193794353
$ws.Range("M2").Value = 17.55809320917071
$ws.Range("O2").Value = 24.74635415821763

$ws.Range("B3").Value = 18.64331542390315
$ws.Range("D3").Value = 10.36993450266866
$ws.Range("E3").Value = 17.45213991047067
$ws.Range("F3").Value = 33.16334309409126
$ws.Range("G3").Value = 3.65538680663969
$ws.Range("J3").Value = 12.30753864399268
$ws.Range("L3").Value = 10.9705724275762
$ws.Range("M3").Value = 17.28350375573791
$ws.Range("O3").Value = 24.86188490544451

$ws.Range("B4").Value = 18.41629147532288
$ws.Range("D4").Value = 10.37198723232048
$ws.Range("E4").Value = 17.39251751273867
$ws.Range("F4").Value = 33.30562224598224
$ws.Range("G4").Value = 3.657030340580299
$ws.Range("J4").Value = 12.27850920928888
$ws.Range("L4").Value = 10.77712349738554
$ws.Range("M4").Value = 17.1131873609417
$ws.Range("O4").Value = 24.94061176770542

$ws.Range("B5").Value = 18.32327893578884
$ws.Range("D5").Value = 10.37307575664721
$ws.Range("E5").Value = 17.36863950231668
$ws.Range("F5").Value = 33.36625986296754
$ws.Range("G5").Value = 3.657720959665581
$ws.Range("J5").Value = 12.26696189043359
$ws.Range("L5").Value = 10.69715973970805
$ws.Range("M5").Value = 17.04341647258398
$ws.Range("O5").Value = 24.97464315729449

$ws.Range("B6").Value = 18.30780733980018
$ws.Range("D6").Value = 10.3732717616587
$ws.Range("E6").Value = 17.36470037019976
$ws.Range("F6").Value = 33.37648891085184
$ws.Range("G6").Value = 3.657836898703992
$ws.Range("J6").Value = 12.26506173890747
$ws.Range("L6").Value = 10.68381626078423
$ws.Range("M6").Value = 17.03181095221206
$ws.Range("O6").Value = 24.98041151234931

$ws.Range("B7").Value = 18.41503895104485
$ws.Range("D7").Value = 10.3720008905314
$ws.Range("E7").Value = 17.39219376748381
$ws.Range("F7").Value = 33.30642927883692
$ws.Range("G7").Value = 3.657039569935149
$ws.Range("J7").Value = 12.27835232488808
$ws.Range("L7").Value = 10.77604953511589
$ws.Range("M7").Value = 17.11224779804347
$ws.Range("O7").Value = 24.94106284592132

$ws.Range("B8").Value = 18.88356680099958
$ws.Range("D8").Value = 10.3686020850393
$ws.Range("E8").Value = 17.51714343547158
$ws.Range("F8").Value = 33.02063176590545
$ws.Range("G8").Value = 3.653704130841105
$ws.Range("J8").Value = 12.33946000057794
$ws.Range("L8").Value = 11.17287254544251
$ws.Range("M8").Value = 17.46381052127926
$ws.Range("O8").Value = 24.78456662476422

$ws.Range("B9").Value = 19.77877542429382
$ws.Range("D9").Value = 10.36979682060082
$ws.Range("E9").Value = 17.77401550249656
$ws.Range("F9").Value = 32.5448406575104
$ws.Range("G9").Value = 3.647816410554611
$ws.Range("J9").Value = 12.4676272693889
$ws.Range("L9").Value = 11.90876919708992
$ws.Range("M9").Value = 18.13676626832919
$ws.Range("O9").Value = 24.53990520446075

$ws.Range("B10").Value = 20.41639586727928
$ws.Range("D10").Value = 10.37539285567331
$ws.Range("E10").Value = 17.96886615589916
$ws.Range("F10").Value = 32.24753165672116
$ws.Range("G10").Value = 3.643884414363683
$ws.Range("J10").Value = 12.56640830399648
$ws.Range("L10").Value = 12.41920680676978
$ws.Range("M10").Value = 18.61749710106931
$ws.Range("O10").Value = 24.39863864442518

$ws.Range("B11").Value = 20.70098654198122
$ws.Range("D11").Value = 10.37894685867074
$ws.Range("E11").Value = 18.05859990466431
$ws.Range("F11").Value = 32.12376589904281
$ws.Range("G11").Value = 3.642180197788356
$ws.Range("J11").Value = 12.61224597811637
$ws.Range("L11").Value = 12.64404370106235
$ws.Range("M11").Value = 18.83249075680226
$ws.Range("O11").Value = 24.34284258784767

$ws.Range("B12").Value = 20.80788124440832
$ws.Range("D12").Value = 10.38043620481798
$ws.Range("E12").Value = 18.09271670030615
$ws.Range("F12").Value = 32.07856208487681
$ws.Range("G12").Value = 3.641546930139191
$ws.Range("J12").Value = 12.62972419118388
$ws.Range("L12").Value = 12.72806539115746
$ws.Range("M12").Value = 18.9133144469028
$ws.Range("O12").Value = 24.32294045203538

$ws.Range("B13").Value = 20.78489974872854
$ws.Range("D13").Value = 10.38010908820513
$ws.Range("E13").Value = 18.08536329798156
$ws.Range("F13").Value = 32.08822334531138
$ws.Range("G13").Value = 3.641682779318429
$ws.Range("J13").Value = 12.62595473440219
$ws.Range("L13").Value = 12.71002048812147
$ws.Range("M13").Value = 18.89593475224708
$ws.Range("O13").Value = 24.32717203745343

$ws.Range("B14").Value = 20.7097988243928
$ws.Range("D14").Value = 10.37906652026674
$ws.Range("E14").Value = 18.06140409704109
$ws.Range("F14").Value = 32.12001354089973
$ws.Range("G14").Value = 3.6421278566997
$ws.Range("J14").Value = 12.61368155965496
$ws.Range("L14").Value = 12.65097896972987
$ws.Range("M14").Value = 18.83915231578235
$ws.Range("O14").Value = 24.34118059292232

$ws.Range("B15").Value = 20.66368106271063
$ws.Range("D15").Value = 10.37844656367125
$ws.Range("E15").Value = 18.04674555779473
$ws.Range("F15").Value = 32.13970297681859
$ws.Range("G15").Value = 3.642402051004079
$ws.Range("J15").Value = 12.60617930120082
$ws.Range("L15").Value = 12.61466690045935
$ws.Range("M15").Value = 18.80429293884705
$ws.Range("O15").Value = 24.34992124419912

$ws.Range("B16").Value = 20.39767959489067
$ws.Range("D16").Value = 10.37518075835766
$ws.Range("E16").Value = 17.96302209741699
$ws.Range("F16").Value = 32.25585214964684
$ws.Range("G16").Value = 3.64399748314011
$ws.Range("J16").Value = 12.56343006667779
$ws.Range("L16").Value = 12.40435972064695
$ws.Range("M16").Value = 18.60336736520885
$ws.Range("O16").Value = 24.4024562074947

$ws.Range("B17").Value = 20.2330322848142
$ws.Range("D17").Value = 10.37343451522687
$ws.Range("E17").Value = 17.91192613309836
$ws.Range("F17").Value = 32.33005551356344
$ws.Range("G17").Value = 3.644997817570952
$ws.Range("J17").Value = 12.53742925062327
$ws.Range("L17").Value = 12.273412248278
$ws.Range("M17").Value = 18.47911749930796
$ws.Range("O17").Value = 24.43685983986389

$ws.Range("B18").Value = 20.13782191455267
$ws.Range("D18").Value = 10.37252512278377
$ws.Range("E18").Value = 17.88264143342343
$ws.Range("F18").Value = 32.37381527002153
$ws.Range("G18").Value = 3.645581137838978
$ws.Range("J18").Value = 12.52255969838488
$ws.Range("L18").Value = 12.19740482597598
$ws.Range("M18").Value = 18.40730797702846
$ws.Range("O18").Value = 24.45744431097749

$ws.Range("B19").Value = 20.10550055365538
$ws.Range("D19").Value = 10.37223357883883
$ws.Range("E19").Value = 17.87274469442913
$ws.Range("F19").Value = 32.38881667255378
$ws.Range("G19").Value = 3.645780008192697
$ws.Range("J19").Value = 12.51754007360836
$ws.Range("L19").Value = 12.17155346303178
$ws.Range("M19").Value = 18.38293720803694
$ws.Range("O19").Value = 24.46455036100711

$ws.Range("B20").Value = 20.25061266458181
$ws.Range("D20").Value = 10.37361058416555
$ws.Range("E20").Value = 17.91735473729764
$ws.Range("F20").Value = 32.32204457798913
$ws.Range("G20").Value = 3.644890507482856
$ws.Range("J20").Value = 12.54018830632833
$ws.Range("L20").Value = 12.28742368788011
$ws.Range("M20").Value = 18.49238017566267
$ws.Range("O20").Value = 24.43311502179967

$ws.Range("B21").Value = 20.73188216707206
$ws.Range("D21").Value = 10.37936886454566
$ws.Range("E21").Value = 18.06843795581628
$ws.Range("F21").Value = 32.11063073985833
$ws.Range("G21").Value = 3.641996799286658
$ws.Range("J21").Value = 12.61728329257549
$ws.Range("L21").Value = 12.66835171930255
$ws.Range("M21").Value = 18.8558471528322
$ws.Range("O21").Value = 24.33703257961316

$ws.Range("B22").Value = 21.04129075037123
$ws.Range("D22").Value = 10.38396807894018
$ws.Range("E22").Value = 18.16796586001453
$ws.Range("F22").Value = 31.98216216192415
$ws.Range("G22").Value = 3.640175988384561
$ws.Range("J22").Value = 12.66836723219557
$ws.Range("L22").Value = 12.91076306575614
$ws.Range("M22").Value = 19.08993115710492
$ws.Range("O22").Value = 24.28139080931547

$ws.Range("B23").Value = 20.87665039218172
$ws.Range("D23").Value = 10.38143740957585
$ws.Range("E23").Value = 18.11478093311067
$ws.Range("F23").Value = 32.04983629787949
$ws.Range("G23").Value = 3.641141369431502
$ws.Range("J23").Value = 12.64104197850449
$ws.Range("L23").Value = 12.78200101399028
$ws.Range("M23").Value = 18.96533160526702
$ws.Range("O23").Value = 24.31043038958479

$ws.Range("B24").Value = 20.24266628885141
$ws.Range("D24").Value = 10.37353068879053
$ws.Range("E24").Value = 17.91490018013216
$ws.Range("F24").Value = 32.3256628986972
$ws.Range("G24").Value = 3.644938996780021
$ws.Range("J24").Value = 12.53894069151267
$ws.Range("L24").Value = 12.28109136645829
$ws.Range("M24").Value = 18.48638528826918
$ws.Range("O24").Value = 24.43480554600242

$ws.Range("B25").Value = 19.53972958362175
$ws.Range("D25").Value = 10.3686395782295
$ws.Range("E25").Value = 17.70336828645773
$ws.Range("F25").Value = 32.66442549419114
$ws.Range("G25").Value = 3.649339737595832
$ws.Range("J25").Value = 12.43210945725478
$ws.Range("L25").Value = 11.71470086287455
$ws.Range("M25").Value = 17.95686161877594
$ws.Range("O25").Value = 24.59937290847238
